$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Runmode column (D) flips from "N" to "Y" for the existing test case rows
# (rows 2 and 3 were already "Y" and are left untouched).
for ($r = 4; $r -le 22; $r++) {
    $ws.Cells.Item($r, 4).Value = "Y"
}

# First copy the formatting of the last existing data row (22) down onto the
# two new rows so the appended rows carry the same cell styles.
$ws.Range("A22:E22").Copy()
$ws.Range("A23:E24").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Two new notification test cases appended below the existing data.
$ws.Cells.Item(23, 1).Value = "Notifications023"
$ws.Cells.Item(23, 2).Value = "OPQA-211"
$ws.Cells.Item(23, 3).Value = "Verify that user is able to view top commenters information in home page"
$ws.Cells.Item(23, 4).Value = "Y"
$ws.Cells.Item(23, 5).Value = "SKIP"

$ws.Cells.Item(24, 1).Value = "Notifications024"
$ws.Cells.Item(24, 2).Value = "OPQA-212"
$ws.Cells.Item(24, 3).Value = "Verify that user is able to view Most viewed documents in home page"
$ws.Cells.Item(24, 4).Value = "Y"
$ws.Cells.Item(24, 5).Value = "SKIP"

# Selection moves from D5 to D4 in the final workbook.
$ws.Range("D4").Select()
